$d = $word.ActiveDocument

$pairs = @(
    @("343×2=", "544×4="),
    @("860×3=", "984×6="),
    @("763×3=", "553×4="),
    @("683×8=", "293×9="),
    @("706×5=", "861×3="),
    @("632×7=", "139×4="),
    @("586×8=", "657×8="),
    @("429×9=", "307×9="),
    @("729×4=", "281×5="),
    @("870×2=", "845×4="),
    @("296×5=", "160×4="),
    @("628×3=", "500×8="),
    @("815×6=", "446×9="),
    @("685×5=", "420×4="),
    @("636×8=", "900×4="),
    @("236×7=", "154×4="),
    @("102×2=", "877×6="),
    @("770×6=", "649×6="),
    @("382×3=", "796×5="),
    @("945×7=", "355×5="),
    @("544×5=", "203×9="),
    @("585×5=", "602×9="),
    @("569×5=", "122×2="),
    @("975×6=", "241×5="),
    @("394×6=", "738×9=")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
